# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.438.16"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").Value = "1.644.11"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "298.89"
$ws.Range("E6").Value = "  -1.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3779"
$ws.Range("E7").Value = "  -0.99%  "

$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("E9").Value = "  -2.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08070"
$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.211"
$ws.Range("E11").Value = "  -3.60%  "

$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.05"
$ws.Range("E13").Value = "  -2.80%  "

$ws.Range("E14").Value = "  -2.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.301"
$ws.Range("E15").Value = "  -2.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001199"
$ws.Range("E16").Value = "  -3.27%  "

$ws.Range("D17").Value = "1.638.37"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06977"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.720"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  -2.27%  "

$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.37"
$ws.Range("E23").Value = "  -3.03%  "

$ws.Range("D24").Value = "23.450.58"
$ws.Range("E24").Value = "  -1.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.478"
$ws.Range("E25").Value = "  -3.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.885"
$ws.Range("E26").Value = "  -6.35%  "

$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.08"
$ws.Range("E28").Value = "  +1.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.203"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.17"
$ws.Range("E30").Value = "  -1.61%  "

$ws.Range("D31").Value = "1.820.02"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.856"
$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.133"
$ws.Range("E33").Value = "  -2.31%  "

$ws.Range("E34").Value = "  -2.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9836"
$ws.Range("E35").Value = "  -9.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02710"
$ws.Range("E36").Value = "  -4.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08733"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("E38").Value = "  -3.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.905"
$ws.Range("E39").Value = "  -4.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06789"
$ws.Range("E40").Value = "  -5.45%  "

$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6854"
$ws.Range("E42").Value = "  -3.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.288"
$ws.Range("E43").Value = "  -4.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.61"
$ws.Range("E44").Value = "  -2.88%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6332"
$ws.Range("E46").Value = "  -3.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.250"
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.901"
$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07714"
$ws.Range("E49").Value = "  -3.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.96"
$ws.Range("E50").Value = "  -1.26%  "

$ws.Range("E51").Value = "  -4.36%  "
